# Auto-generated Excel COM-interop edit script
# Applies refreshed currentAveragePrice / Leve price / profit figures
# from an external market-data snapshot onto the Leve profit sheets.
#
# Each block below:
#   - selects the worksheet by name (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR)
#   - verifies the row via the stable "Leve Item ID" in column G
#   - writes the updated H/I/J/K/L/M/N values (or clears ones that the
#     snapshot no longer has data for, matching upstream behavior)
$wb = $excel.ActiveWorkbook

# ALC row 74 (Leve Item ID = 5507)
$ws = $wb.Worksheets.Item("ALC")
if ($ws.Cells.Item(74, 7).Value2 -ne 5507) { throw "Unexpected G74 on ALC: expected 5507" }
$ws.Cells.Item(74, 8).Value2 = 3767.1667  # H74: 3704 -> 3767.1667
$ws.Cells.Item(74, 9).Value2 = 3500  # I74: 3350.75 -> 3500
$ws.Cells.Item(74, 10).Value2 = 3820.6  # J74: 4410.5 -> 3820.6
$ws.Cells.Item(74, 11).Value2 = 3500  # K74: 3350.75 -> 3500
$ws.Cells.Item(74, 12).Value2 = 3820.6  # L74: 4410.5 -> 3820.6
$ws.Cells.Item(74, 13).Value2 = -2564  # M74: -2414.75 -> -2564
$ws.Cells.Item(74, 14).Value2 = -5692.6  # N74: -6282.5 -> -5692.6

# ALC row 77 (Leve Item ID = 5507)
$ws = $wb.Worksheets.Item("ALC")
if ($ws.Cells.Item(77, 7).Value2 -ne 5507) { throw "Unexpected G77 on ALC: expected 5507" }
$ws.Cells.Item(77, 8).Value2 = 3767.1667  # H77: 3704 -> 3767.1667
$ws.Cells.Item(77, 9).Value2 = 3500  # I77: 3350.75 -> 3500
$ws.Cells.Item(77, 10).Value2 = 3820.6  # J77: 4410.5 -> 3820.6
$ws.Cells.Item(77, 11).Value2 = 17500  # K77: 16753.75 -> 17500
$ws.Cells.Item(77, 12).Value2 = 19103  # L77: 22052.5 -> 19103
$ws.Cells.Item(77, 13).Value2 = -12820  # M77: -12073.75 -> -12820
$ws.Cells.Item(77, 14).Value2 = -28463  # N77: -31412.5 -> -28463

# ALC row 129 (Leve Item ID = 36115)
$ws = $wb.Worksheets.Item("ALC")
if ($ws.Cells.Item(129, 7).Value2 -ne 36115) { throw "Unexpected G129 on ALC: expected 36115" }
$ws.Cells.Item(129, 8).Value2 = 1108.0377  # H129: 1158 -> 1108.0377
$ws.Cells.Item(129, 9).Value2 = 296  # I129: 314.66666 -> 296
$ws.Cells.Item(129, 10).Value2 = 1174.3265  # J129: 1222.8718 -> 1174.3265
$ws.Cells.Item(129, 11).Value2 = 888  # K129: 943.9999799999999 -> 888
$ws.Cells.Item(129, 12).Value2 = 3522.979499999999  # L129: 3668.6154 -> 3522.979499999999
$ws.Cells.Item(129, 13).Value2 = 4112  # M129: 4056.00002 -> 4112
$ws.Cells.Item(129, 14).Value2 = -13522.9795  # N129: -13668.6154 -> -13522.9795

# ALC row 131 (Leve Item ID = 36108)
$ws = $wb.Worksheets.Item("ALC")
if ($ws.Cells.Item(131, 7).Value2 -ne 36108) { throw "Unexpected G131 on ALC: expected 36108" }
$ws.Cells.Item(131, 8).Value2 = 4112.143  # H131: 3603.3333 -> 4112.143
$ws.Cells.Item(131, 9).Value2 = 1695  # I131: 1933.75 -> 1695
$ws.Cells.Item(131, 10).Value2 = 6309.5454  # J131: 4630.769 -> 6309.5454
$ws.Cells.Item(131, 11).Value2 = 5085  # K131: 5801.25 -> 5085
$ws.Cells.Item(131, 12).Value2 = 18928.6362  # L131: 13892.307 -> 18928.6362
$ws.Cells.Item(131, 13).Value2 = -45  # M131: -761.25 -> -45
$ws.Cells.Item(131, 14).Value2 = -29008.6362  # N131: -23972.307 -> -29008.6362

# ALC row 137 (Leve Item ID = 44013)
$ws = $wb.Worksheets.Item("ALC")
if ($ws.Cells.Item(137, 7).Value2 -ne 44013) { throw "Unexpected G137 on ALC: expected 44013" }
$ws.Cells.Item(137, 8).Value2 = 1525.3103  # H137: 1736.8334 -> 1525.3103
$ws.Cells.Item(137, 9).Value2 = 1004.8  # I137: 1144.6 -> 1004.8
$ws.Cells.Item(137, 10).Value2 = 2083  # J137: 2159.8572 -> 2083
$ws.Cells.Item(137, 11).Value2 = 3014.4  # K137: 3433.8 -> 3014.4
$ws.Cells.Item(137, 12).Value2 = 6249  # L137: 6479.571599999999 -> 6249
$ws.Cells.Item(137, 13).Value2 = -464.3999999999996  # M137: -883.7999999999997 -> -464.3999999999996
$ws.Cells.Item(137, 14).Value2 = -11349  # N137: -11579.5716 -> -11349

# ALC row 139 (Leve Item ID = 42306)
$ws = $wb.Worksheets.Item("ALC")
if ($ws.Cells.Item(139, 7).Value2 -ne 42306) { throw "Unexpected G139 on ALC: expected 42306" }
$ws.Cells.Item(139, 8).Value2 = 77875  # H139: 78040 -> 77875
$ws.Cells.Item(139, 10).Value2 = 77875  # J139: 78040 -> 77875
$ws.Cells.Item(139, 12).Value2 = 77875  # L139: 78040 -> 77875
$ws.Cells.Item(139, 14).Value2 = -88155  # N139: -88320 -> -88155

# ARM row 61 (Leve Item ID = 43999)
$ws = $wb.Worksheets.Item("ARM")
if ($ws.Cells.Item(61, 7).Value2 -ne 43999) { throw "Unexpected G61 on ARM: expected 43999" }
$ws.Cells.Item(61, 8).Value2 = 2052.4722  # H61: 2254.484 -> 2052.4722
$ws.Cells.Item(61, 9).Value2 = 2086.7188  # I61: 2325 -> 2086.7188
$ws.Cells.Item(61, 11).Value2 = 2086.7188  # K61: 2325 -> 2086.7188
$ws.Cells.Item(61, 13).Value2 = -1874.7188  # M61: -2113 -> -1874.7188

# ARM row 74 (Leve Item ID = 44000)
$ws = $wb.Worksheets.Item("ARM")
if ($ws.Cells.Item(74, 7).Value2 -ne 44000) { throw "Unexpected G74 on ARM: expected 44000" }
$ws.Cells.Item(74, 8).Value2 = 971  # H74: 1039.4286 -> 971
$ws.Cells.Item(74, 9).Value2 = 862.4  # I74: 960.3333 -> 862.4
$ws.Cells.Item(74, 11).Value2 = 862.4  # K74: 960.3333 -> 862.4
$ws.Cells.Item(74, 13).Value2 = 11.60000000000002  # M74: -86.33330000000001 -> 11.60000000000002

# ARM row 77 (Leve Item ID = 44000)
$ws = $wb.Worksheets.Item("ARM")
if ($ws.Cells.Item(77, 7).Value2 -ne 44000) { throw "Unexpected G77 on ARM: expected 44000" }
$ws.Cells.Item(77, 8).Value2 = 971  # H77: 1039.4286 -> 971
$ws.Cells.Item(77, 9).Value2 = 862.4  # I77: 960.3333 -> 862.4
$ws.Cells.Item(77, 11).Value2 = 4312  # K77: 4801.6665 -> 4312
$ws.Cells.Item(77, 13).Value2 = 56  # M77: -433.6665000000003 -> 56

# ARM row 97 (Leve Item ID = 19941)
$ws = $wb.Worksheets.Item("ARM")
if ($ws.Cells.Item(97, 7).Value2 -ne 19941) { throw "Unexpected G97 on ARM: expected 19941" }
$ws.Cells.Item(97, 8).Value2 = 1836.1428  # H97: 1985.95 -> 1836.1428
$ws.Cells.Item(97, 9).Value2 = 1503.2778  # I97: 1679.9395 -> 1503.2778
$ws.Cells.Item(97, 10).Value2 = 3833.3333  # J97: 3428.5715 -> 3833.3333
$ws.Cells.Item(97, 11).Value2 = 1503.2778  # K97: 1679.9395 -> 1503.2778
$ws.Cells.Item(97, 12).Value2 = 3833.3333  # L97: 3428.5715 -> 3833.3333
$ws.Cells.Item(97, 13).Value2 = -1007.2778  # M97: -1183.9395 -> -1007.2778
$ws.Cells.Item(97, 14).Value2 = -4825.3333  # N97: -4420.5715 -> -4825.3333

# ARM row 133 (Leve Item ID = 41857)
$ws = $wb.Worksheets.Item("ARM")
if ($ws.Cells.Item(133, 7).Value2 -ne 41857) { throw "Unexpected G133 on ARM: expected 41857" }
$ws.Cells.Item(133, 8).Value2 = 76026.10000000001  # H133: 70558.53999999999 -> 76026.10000000001
$ws.Cells.Item(133, 10).Value2 = 76026.10000000001  # J133: 70558.53999999999 -> 76026.10000000001
$ws.Cells.Item(133, 12).Value2 = 76026.10000000001  # L133: 70558.53999999999 -> 76026.10000000001
$ws.Cells.Item(133, 14).Value2 = -81086.10000000001  # N133: -75618.53999999999 -> -81086.10000000001

# ARM row 134 (Leve Item ID = 42019)
$ws = $wb.Worksheets.Item("ARM")
if ($ws.Cells.Item(134, 7).Value2 -ne 42019) { throw "Unexpected G134 on ARM: expected 42019" }
$ws.Cells.Item(134, 8).Value2 = 0  # H134: 39990 -> 0
$ws.Cells.Item(134, 10).Value2 = 0  # J134: 39990 -> 0
$ws.Cells.Item(134, 12).Value2 = 0  # L134: 39990 -> 0
$ws.Cells.Item(134, 14).ClearContents()  # N134: -50130 -> (removed)

# ARM row 135 (Leve Item ID = 42016)
$ws = $wb.Worksheets.Item("ARM")
if ($ws.Cells.Item(135, 7).Value2 -ne 42016) { throw "Unexpected G135 on ARM: expected 42016" }
$ws.Cells.Item(135, 8).Value2 = 45952.668  # H135: 45631.89 -> 45952.668
$ws.Cells.Item(135, 10).Value2 = 45952.668  # J135: 45631.89 -> 45952.668
$ws.Cells.Item(135, 12).Value2 = 45952.668  # L135: 45631.89 -> 45952.668
$ws.Cells.Item(135, 14).Value2 = -56092.668  # N135: -55771.89 -> -56092.668

# ARM row 136 (Leve Item ID = 43999)
$ws = $wb.Worksheets.Item("ARM")
if ($ws.Cells.Item(136, 7).Value2 -ne 43999) { throw "Unexpected G136 on ARM: expected 43999" }
$ws.Cells.Item(136, 8).Value2 = 2052.4722  # H136: 2254.484 -> 2052.4722
$ws.Cells.Item(136, 9).Value2 = 2086.7188  # I136: 2325 -> 2086.7188
$ws.Cells.Item(136, 11).Value2 = 6260.1564  # K136: 6975 -> 6260.1564
$ws.Cells.Item(136, 13).Value2 = -3710.1564  # M136: -4425 -> -3710.1564

# ARM row 139 (Leve Item ID = 42321)
$ws = $wb.Worksheets.Item("ARM")
if ($ws.Cells.Item(139, 7).Value2 -ne 42321) { throw "Unexpected G139 on ARM: expected 42321" }
$ws.Cells.Item(139, 8).Value2 = 55000  # H139: 43857.5 -> 55000
$ws.Cells.Item(139, 10).Value2 = 55000  # J139: 43857.5 -> 55000
$ws.Cells.Item(139, 12).Value2 = 55000  # L139: 43857.5 -> 55000
$ws.Cells.Item(139, 14).Value2 = -65280  # N139: -54137.5 -> -65280

# ARM row 140 (Leve Item ID = 42496)
$ws = $wb.Worksheets.Item("ARM")
if ($ws.Cells.Item(140, 7).Value2 -ne 42496) { throw "Unexpected G140 on ARM: expected 42496" }
$ws.Cells.Item(140, 8).Value2 = 103032.9  # H140: 108381 -> 103032.9
$ws.Cells.Item(140, 10).Value2 = 103032.9  # J140: 108381 -> 103032.9
$ws.Cells.Item(140, 12).Value2 = 103032.9  # L140: 108381 -> 103032.9
$ws.Cells.Item(140, 14).Value2 = -113392.9  # N140: -118741 -> -113392.9

# BSM row 132 (Leve Item ID = 41855)
$ws = $wb.Worksheets.Item("BSM")
if ($ws.Cells.Item(132, 7).Value2 -ne 41855) { throw "Unexpected G132 on BSM: expected 41855" }
$ws.Cells.Item(132, 8).Value2 = 0  # H132: 36200 -> 0
$ws.Cells.Item(132, 10).Value2 = 0  # J132: 36200 -> 0
$ws.Cells.Item(132, 12).Value2 = 0  # L132: 36200 -> 0
$ws.Cells.Item(132, 14).ClearContents()  # N132: -46320 -> (removed)

# BSM row 135 (Leve Item ID = 41992)
$ws = $wb.Worksheets.Item("BSM")
if ($ws.Cells.Item(135, 7).Value2 -ne 41992) { throw "Unexpected G135 on BSM: expected 41992" }
$ws.Cells.Item(135, 8).Value2 = 50000  # H135: 45862.223 -> 50000
$ws.Cells.Item(135, 10).Value2 = 50000  # J135: 45862.223 -> 50000
$ws.Cells.Item(135, 12).Value2 = 50000  # L135: 45862.223 -> 50000
$ws.Cells.Item(135, 14).Value2 = -60140  # N135: -56002.223 -> -60140

# BSM row 137 (Leve Item ID = 42153)
$ws = $wb.Worksheets.Item("BSM")
if ($ws.Cells.Item(137, 7).Value2 -ne 42153) { throw "Unexpected G137 on BSM: expected 42153" }
$ws.Cells.Item(137, 8).Value2 = 85000  # H137: 0 -> 85000
$ws.Cells.Item(137, 10).Value2 = 85000  # J137: 0 -> 85000
$ws.Cells.Item(137, 12).Value2 = 85000  # L137: 0 -> 85000
$ws.Cells.Item(137, 14).Value2 = -95200  # N137: None -> -95200

# CRP row 31 (Leve Item ID = 44023)
$ws = $wb.Worksheets.Item("CRP")
if ($ws.Cells.Item(31, 7).Value2 -ne 44023) { throw "Unexpected G31 on CRP: expected 44023" }
$ws.Cells.Item(31, 8).Value2 = 11610.523  # H31: 12528.167 -> 11610.523
$ws.Cells.Item(31, 9).Value2 = 4501.1665  # I31: 5291.263 -> 4501.1665
$ws.Cells.Item(31, 10).Value2 = 15772.098  # J31: 15881.854 -> 15772.098
$ws.Cells.Item(31, 11).Value2 = 4501.1665  # K31: 5291.263 -> 4501.1665
$ws.Cells.Item(31, 12).Value2 = 15772.098  # L31: 15881.854 -> 15772.098
$ws.Cells.Item(31, 13).Value2 = -4206.1665  # M31: -4996.263 -> -4206.1665
$ws.Cells.Item(31, 14).Value2 = -16362.098  # N31: -16471.854 -> -16362.098

# CRP row 34 (Leve Item ID = 44023)
$ws = $wb.Worksheets.Item("CRP")
if ($ws.Cells.Item(34, 7).Value2 -ne 44023) { throw "Unexpected G34 on CRP: expected 44023" }
$ws.Cells.Item(34, 8).Value2 = 11610.523  # H34: 12528.167 -> 11610.523
$ws.Cells.Item(34, 9).Value2 = 4501.1665  # I34: 5291.263 -> 4501.1665
$ws.Cells.Item(34, 10).Value2 = 15772.098  # J34: 15881.854 -> 15772.098
$ws.Cells.Item(34, 11).Value2 = 4501.1665  # K34: 5291.263 -> 4501.1665
$ws.Cells.Item(34, 12).Value2 = 15772.098  # L34: 15881.854 -> 15772.098
$ws.Cells.Item(34, 13).Value2 = -4299.1665  # M34: -5089.263 -> -4299.1665
$ws.Cells.Item(34, 14).Value2 = -16176.098  # N34: -16285.854 -> -16176.098

# CRP row 130 (Leve Item ID = 34689)
$ws = $wb.Worksheets.Item("CRP")
if ($ws.Cells.Item(130, 7).Value2 -ne 34689) { throw "Unexpected G130 on CRP: expected 34689" }
$ws.Cells.Item(130, 8).Value2 = 0  # H130: 52593.332 -> 0
$ws.Cells.Item(130, 10).Value2 = 0  # J130: 52593.332 -> 0
$ws.Cells.Item(130, 12).Value2 = 0  # L130: 52593.332 -> 0
$ws.Cells.Item(130, 14).ClearContents()  # N130: -62633.332 -> (removed)

# CRP row 132 (Leve Item ID = 44019)
$ws = $wb.Worksheets.Item("CRP")
if ($ws.Cells.Item(132, 7).Value2 -ne 44019) { throw "Unexpected G132 on CRP: expected 44019" }
$ws.Cells.Item(132, 8).Value2 = 2479.6924  # H132: 3180.2222 -> 2479.6924
$ws.Cells.Item(132, 9).Value2 = 1223.8  # I132: 1437.3334 -> 1223.8
$ws.Cells.Item(132, 11).Value2 = 3671.4  # K132: 4312.0002 -> 3671.4
$ws.Cells.Item(132, 13).Value2 = -1141.4  # M132: -1782.0002 -> -1141.4

# CRP row 133 (Leve Item ID = 43328)
$ws = $wb.Worksheets.Item("CRP")
if ($ws.Cells.Item(133, 7).Value2 -ne 43328) { throw "Unexpected G133 on CRP: expected 43328" }
$ws.Cells.Item(133, 8).Value2 = 0  # H133: 35000 -> 0
$ws.Cells.Item(133, 10).Value2 = 0  # J133: 35000 -> 0
$ws.Cells.Item(133, 12).Value2 = 0  # L133: 35000 -> 0
$ws.Cells.Item(133, 14).ClearContents()  # N133: -40060 -> (removed)

# CRP row 135 (Leve Item ID = 42008)
$ws = $wb.Worksheets.Item("CRP")
if ($ws.Cells.Item(135, 7).Value2 -ne 42008) { throw "Unexpected G135 on CRP: expected 42008" }
$ws.Cells.Item(135, 8).Value2 = 0  # H135: 46225 -> 0
$ws.Cells.Item(135, 10).Value2 = 0  # J135: 46225 -> 0
$ws.Cells.Item(135, 12).Value2 = 0  # L135: 46225 -> 0
$ws.Cells.Item(135, 14).ClearContents()  # N135: -56365 -> (removed)

# CRP row 138 (Leve Item ID = 42302)
$ws = $wb.Worksheets.Item("CRP")
if ($ws.Cells.Item(138, 7).Value2 -ne 42302) { throw "Unexpected G138 on CRP: expected 42302" }
$ws.Cells.Item(138, 8).Value2 = 0  # H138: 49956 -> 0
$ws.Cells.Item(138, 10).Value2 = 0  # J138: 49956 -> 0
$ws.Cells.Item(138, 12).Value2 = 0  # L138: 49956 -> 0
$ws.Cells.Item(138, 14).ClearContents()  # N138: -60236 -> (removed)

# CRP row 140 (Leve Item ID = 42455)
$ws = $wb.Worksheets.Item("CRP")
if ($ws.Cells.Item(140, 7).Value2 -ne 42455) { throw "Unexpected G140 on CRP: expected 42455" }
$ws.Cells.Item(140, 8).Value2 = 90000  # H140: 65000.5 -> 90000
$ws.Cells.Item(140, 10).Value2 = 90000  # J140: 65000.5 -> 90000
$ws.Cells.Item(140, 12).Value2 = 90000  # L140: 65000.5 -> 90000
$ws.Cells.Item(140, 14).Value2 = -100360  # N140: -75360.5 -> -100360

# GSM row 122 (Leve Item ID = 36182)
$ws = $wb.Worksheets.Item("GSM")
if ($ws.Cells.Item(122, 7).Value2 -ne 36182) { throw "Unexpected G122 on GSM: expected 36182" }
$ws.Cells.Item(122, 8).Value2 = 3371.75  # H122: 2053.5715 -> 3371.75
$ws.Cells.Item(122, 9).Value2 = 2829  # I122: 1562.5 -> 2829
$ws.Cells.Item(122, 11).Value2 = 8487  # K122: 4687.5 -> 8487
$ws.Cells.Item(122, 13).Value2 = -6037  # M122: -2237.5 -> -6037

# GSM row 132 (Leve Item ID = 44008)
$ws = $wb.Worksheets.Item("GSM")
if ($ws.Cells.Item(132, 7).Value2 -ne 44008) { throw "Unexpected G132 on GSM: expected 44008" }
$ws.Cells.Item(132, 8).Value2 = 2727.0952  # H132: 3026 -> 2727.0952
$ws.Cells.Item(132, 9).Value2 = 1872.8  # I132: 2107.5833 -> 1872.8
$ws.Cells.Item(132, 11).Value2 = 5618.4  # K132: 6322.749899999999 -> 5618.4
$ws.Cells.Item(132, 13).Value2 = -3088.4  # M132: -3792.749899999999 -> -3088.4

# GSM row 133 (Leve Item ID = 41854)
$ws = $wb.Worksheets.Item("GSM")
if ($ws.Cells.Item(133, 7).Value2 -ne 41854) { throw "Unexpected G133 on GSM: expected 41854" }
$ws.Cells.Item(133, 8).Value2 = 53900  # H133: 54000 -> 53900
$ws.Cells.Item(133, 10).Value2 = 53900  # J133: 54000 -> 53900
$ws.Cells.Item(133, 12).Value2 = 53900  # L133: 54000 -> 53900
$ws.Cells.Item(133, 14).Value2 = -64020  # N133: -64120 -> -64020

# GSM row 135 (Leve Item ID = 42006)
$ws = $wb.Worksheets.Item("GSM")
if ($ws.Cells.Item(135, 7).Value2 -ne 42006) { throw "Unexpected G135 on GSM: expected 42006" }
$ws.Cells.Item(135, 8).Value2 = 0  # H135: 47990 -> 0
$ws.Cells.Item(135, 10).Value2 = 0  # J135: 47990 -> 0
$ws.Cells.Item(135, 12).Value2 = 0  # L135: 47990 -> 0
$ws.Cells.Item(135, 14).ClearContents()  # N135: -58130 -> (removed)

# GSM row 138 (Leve Item ID = 42325)
$ws = $wb.Worksheets.Item("GSM")
if ($ws.Cells.Item(138, 7).Value2 -ne 42325) { throw "Unexpected G138 on GSM: expected 42325" }
$ws.Cells.Item(138, 8).Value2 = 69480  # H138: 69450 -> 69480
$ws.Cells.Item(138, 10).Value2 = 69480  # J138: 69450 -> 69480
$ws.Cells.Item(138, 12).Value2 = 69480  # L138: 69450 -> 69480
$ws.Cells.Item(138, 14).Value2 = -79760  # N138: -79730 -> -79760

# GSM row 139 (Leve Item ID = 42373)
$ws = $wb.Worksheets.Item("GSM")
if ($ws.Cells.Item(139, 7).Value2 -ne 42373) { throw "Unexpected G139 on GSM: expected 42373" }
$ws.Cells.Item(139, 8).Value2 = 53427.8  # H139: 52564.832 -> 53427.8
$ws.Cells.Item(139, 10).Value2 = 53427.8  # J139: 52564.832 -> 53427.8
$ws.Cells.Item(139, 12).Value2 = 53427.8  # L139: 52564.832 -> 53427.8
$ws.Cells.Item(139, 14).Value2 = -63707.8  # N139: -62844.832 -> -63707.8

# GSM row 140 (Leve Item ID = 42458)
$ws = $wb.Worksheets.Item("GSM")
if ($ws.Cells.Item(140, 7).Value2 -ne 42458) { throw "Unexpected G140 on GSM: expected 42458" }
$ws.Cells.Item(140, 8).Value2 = 0  # H140: 75716.336 -> 0
$ws.Cells.Item(140, 10).Value2 = 0  # J140: 75716.336 -> 0
$ws.Cells.Item(140, 12).Value2 = 0  # L140: 75716.336 -> 0
$ws.Cells.Item(140, 14).ClearContents()  # N140: -86076.336 -> (removed)

# GSM row 141 (Leve Item ID = 42504)
$ws = $wb.Worksheets.Item("GSM")
if ($ws.Cells.Item(141, 7).Value2 -ne 42504) { throw "Unexpected G141 on GSM: expected 42504" }
$ws.Cells.Item(141, 8).Value2 = 47949.332  # H141: 45412 -> 47949.332
$ws.Cells.Item(141, 10).Value2 = 47949.332  # J141: 45412 -> 47949.332
$ws.Cells.Item(141, 12).Value2 = 47949.332  # L141: 45412 -> 47949.332
$ws.Cells.Item(141, 14).Value2 = -58309.332  # N141: -55772 -> -58309.332

# LTW row 50 (Leve Item ID = 3426)
$ws = $wb.Worksheets.Item("LTW")
if ($ws.Cells.Item(50, 7).Value2 -ne 3426) { throw "Unexpected G50 on LTW: expected 3426" }
$ws.Cells.Item(50, 8).Value2 = 6625  # H50: 7000 -> 6625
$ws.Cells.Item(50, 10).Value2 = 6625  # J50: 7000 -> 6625
$ws.Cells.Item(50, 12).Value2 = 6625  # L50: 7000 -> 6625
$ws.Cells.Item(50, 14).Value2 = -7899  # N50: -8274 -> -7899

# LTW row 132 (Leve Item ID = 44058)
$ws = $wb.Worksheets.Item("LTW")
if ($ws.Cells.Item(132, 7).Value2 -ne 44058) { throw "Unexpected G132 on LTW: expected 44058" }
$ws.Cells.Item(132, 8).Value2 = 5042.7896  # H132: 3450.9429 -> 5042.7896
$ws.Cells.Item(132, 9).Value2 = 5059.647  # I132: 3675.6072 -> 5059.647
$ws.Cells.Item(132, 10).Value2 = 4899.5  # J132: 2552.2856 -> 4899.5
$ws.Cells.Item(132, 11).Value2 = 15178.941  # K132: 11026.8216 -> 15178.941
$ws.Cells.Item(132, 12).Value2 = 14698.5  # L132: 7656.8568 -> 14698.5
$ws.Cells.Item(132, 13).Value2 = -12648.941  # M132: -8496.821599999999 -> -12648.941
$ws.Cells.Item(132, 14).Value2 = -19758.5  # N132: -12716.8568 -> -19758.5

# LTW row 139 (Leve Item ID = 43310)
$ws = $wb.Worksheets.Item("LTW")
if ($ws.Cells.Item(139, 7).Value2 -ne 43310) { throw "Unexpected G139 on LTW: expected 43310" }
$ws.Cells.Item(139, 8).Value2 = 64540  # H139: 60850 -> 64540
$ws.Cells.Item(139, 10).Value2 = 79425  # J139: 79466.664 -> 79425
$ws.Cells.Item(139, 12).Value2 = 79425  # L139: 79466.664 -> 79425
$ws.Cells.Item(139, 14).Value2 = -89705  # N139: -89746.664 -> -89705

# WVR row 136 (Leve Item ID = 44031)
$ws = $wb.Worksheets.Item("WVR")
if ($ws.Cells.Item(136, 7).Value2 -ne 44031) { throw "Unexpected G136 on WVR: expected 44031" }
$ws.Cells.Item(136, 8).Value2 = 1150.9375  # H136: 986.6389 -> 1150.9375
$ws.Cells.Item(136, 9).Value2 = 731.7646999999999  # I136: 638.2857 -> 731.7646999999999
$ws.Cells.Item(136, 10).Value2 = 1626  # J136: 1474.3334 -> 1626
$ws.Cells.Item(136, 11).Value2 = 2195.2941  # K136: 1914.8571 -> 2195.2941
$ws.Cells.Item(136, 12).Value2 = 4878  # L136: 4423.0002 -> 4878
$ws.Cells.Item(136, 13).Value2 = 354.7058999999999  # M136: 635.1428999999998 -> 354.7058999999999
$ws.Cells.Item(136, 14).Value2 = -9978  # N136: -9523.0002 -> -9978

# WVR row 138 (Leve Item ID = 42347)
$ws = $wb.Worksheets.Item("WVR")
if ($ws.Cells.Item(138, 7).Value2 -ne 42347) { throw "Unexpected G138 on WVR: expected 42347" }
$ws.Cells.Item(138, 8).Value2 = 76440  # H138: 75016.664 -> 76440
$ws.Cells.Item(138, 10).Value2 = 86800  # J138: 83020 -> 86800
$ws.Cells.Item(138, 12).Value2 = 86800  # L138: 83020 -> 86800
$ws.Cells.Item(138, 14).Value2 = -97080  # N138: -93300 -> -97080

Write-Host "Done applying Leve profit updates."